# Slide 4: restyle the "main -> fact(n)" return arrow to match the removed
# "fact(n) -> fact(n-1)" return arrow (flip it horizontally and recolor it
# to accent6), then delete the now-redundant "fact(n-1) call" group and its
# associated return arrow / label textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Shape id=7 "箭头: 右 6" -- the return-to-main arrow.
$returnArrow = $s.Shapes.Item(2)
$returnArrow.HorizontalFlip = -1
$returnArrow.Fill.ForeColor.ObjectThemeColor = 10

# Remove shape id=12 "组合 11" (fact(n-1) call group: rounded rect + "fact(n-1)" label).
$s.Shapes.Item("组合 11").Delete()

# Remove shape id=15 "箭头: 右 14" (the old flipped accent6 return arrow).
$s.Shapes.Item("箭头: 右 14").Delete()

# Remove shape id=16 "文本框 15" (the "ra_n-1: / return to fact(n)" label).
$s.Shapes.Item("文本框 15").Delete()
